$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 459.86517
$ws.Range("J17").Value = 459.86517
$ws.Range("L17").Value = 1379.59551
$ws.Range("N17").Value = -1715.59551
$ws.Range("H40").Value = 1733.25
$ws.Range("I40").Value = 1542.7142
$ws.Range("K40").Value = 1542.7142
$ws.Range("M40").Value = -1367.7142
$ws.Range("H129").Value = 1177.3182
$ws.Range("J129").Value = 1390.8823
$ws.Range("L129").Value = 4172.6469
$ws.Range("N129").Value = -14172.6469
$ws.Range("H137").Value = 1407.3235
$ws.Range("I137").Value = 1417.0769
$ws.Range("J137").Value = 1375.625
$ws.Range("K137").Value = 4251.2307
$ws.Range("L137").Value = 4126.875
$ws.Range("M137").Value = -1701.2307
$ws.Range("N137").Value = -9226.875
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
$ws.Range("H141").Value = 4138.3687
$ws.Range("I141").Value = 1510
$ws.Range("J141").Value = 34802.668
$ws.Range("K141").Value = 4530
$ws.Range("L141").Value = 104408.004
$ws.Range("M141").Value = 650
$ws.Range("N141").Value = -114768.004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1364.4445
$ws.Range("I2").Value = 1047.1428
$ws.Range("J2").Value = 2475
$ws.Range("K2").Value = 1047.1428
$ws.Range("L2").Value = 2475
$ws.Range("M2").Value = -934.1428000000001
$ws.Range("N2").Value = -2701
$ws.Range("H32").Value = 16545.209
$ws.Range("I32").Value = 17554.342
$ws.Range("J32").Value = 10634.571
$ws.Range("K32").Value = 17554.342
$ws.Range("L32").Value = 10634.571
$ws.Range("M32").Value = -17267.342
$ws.Range("N32").Value = -11208.571
$ws.Range("H45").Value = 1024
$ws.Range("I45").Value = 983.17645
$ws.Range("J45").Value = 1197.5
$ws.Range("K45").Value = 983.17645
$ws.Range("L45").Value = 1197.5
$ws.Range("M45").Value = -606.17645
$ws.Range("N45").Value = -1951.5
$ws.Range("H63").Value = 3167.6667
$ws.Range("I63").Value = 2601.2
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 2601.2
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -1915.2
$ws.Range("N63").Value = -7372
$ws.Range("H66").Value = 3167.6667
$ws.Range("I66").Value = 2601.2
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 13006
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -9574
$ws.Range("N66").Value = -36864
$ws.Range("H74").Value = 873.5952
$ws.Range("I74").Value = 740.2857
$ws.Range("J74").Value = 1140.2142
$ws.Range("K74").Value = 740.2857
$ws.Range("L74").Value = 1140.2142
$ws.Range("M74").Value = 133.7143
$ws.Range("N74").Value = -2888.2142
$ws.Range("H77").Value = 873.5952
$ws.Range("I77").Value = 740.2857
$ws.Range("J77").Value = 1140.2142
$ws.Range("K77").Value = 3701.4285
$ws.Range("L77").Value = 5701.071
$ws.Range("M77").Value = 666.5715
$ws.Range("N77").Value = -14437.071
$ws.Range("H116").Value = 1364.4445
$ws.Range("I116").Value = 1047.1428
$ws.Range("J116").Value = 2475
$ws.Range("K116").Value = 1047.1428
$ws.Range("L116").Value = 2475
$ws.Range("M116").Value = 1246.8572
$ws.Range("N116").Value = -7063
$ws.Range("H122").Value = 2506.7932
$ws.Range("I122").Value = 2526.5
$ws.Range("J122").Value = 2444.8572
$ws.Range("K122").Value = 7579.5
$ws.Range("L122").Value = 7334.571599999999
$ws.Range("M122").Value = -5129.5
$ws.Range("N122").Value = -12234.5716
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1364.4445
$ws.Range("I3").Value = 1047.1428
$ws.Range("J3").Value = 2475
$ws.Range("K3").Value = 1047.1428
$ws.Range("L3").Value = 2475
$ws.Range("M3").Value = -933.1428000000001
$ws.Range("N3").Value = -2703
$ws.Range("H94").Value = 1875.375
$ws.Range("I94").Value = 1090.75
$ws.Range("J94").Value = 2660
$ws.Range("K94").Value = 1090.75
$ws.Range("L94").Value = 2660
$ws.Range("M94").Value = -639.75
$ws.Range("N94").Value = -3562
$ws.Range("H99").Value = 1471.4615
$ws.Range("I99").Value = 913
$ws.Range("K99").Value = 913
$ws.Range("M99").Value = 585
$ws.Range("H105").Value = 2580.3333
$ws.Range("I105").Value = 2628.9285
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 2628.9285
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = -881.9285
$ws.Range("N105").Value = -5394
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2392.6924
$ws.Range("I55").Value = 2666.6667
$ws.Range("J55").Value = 2310.5
$ws.Range("K55").Value = 8000.000100000001
$ws.Range("L55").Value = 6931.5
$ws.Range("M55").Value = -7823.000100000001
$ws.Range("N55").Value = -7285.5
$ws.Range("H63").Value = 303087.44
$ws.Range("I63").Value = 702004
$ws.Range("J63").Value = 3900
$ws.Range("K63").Value = 2106012
$ws.Range("L63").Value = 11700
$ws.Range("M63").Value = -2105263
$ws.Range("N63").Value = -13198
$ws.Range("H66").Value = 303087.44
$ws.Range("I66").Value = 702004
$ws.Range("J66").Value = 3900
$ws.Range("K66").Value = 6318036
$ws.Range("L66").Value = 35100
$ws.Range("M66").Value = -6314292
$ws.Range("N66").Value = -42588
$ws.Range("H122").Value = 862
$ws.Range("I122").Value = 486.25
$ws.Range("J122").Value = 1291.4286
$ws.Range("K122").Value = 4376.25
$ws.Range("L122").Value = 11622.8574
$ws.Range("M122").Value = -1926.25
$ws.Range("N122").Value = -16522.8574
$ws.Range("H129").Value = 6250887
$ws.Range("I129").Value = 915
$ws.Range("J129").Value = 8334211
$ws.Range("K129").Value = 2745
$ws.Range("L129").Value = 25002633
$ws.Range("M129").Value = 2255
$ws.Range("N129").Value = -25012633
$ws.Range("H131").Value = 938.35
$ws.Range("J131").Value = 965.5263
$ws.Range("L131").Value = 2896.5789
$ws.Range("N131").Value = -12976.5789
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3064.8948
$ws.Range("I126").Value = 3083.3333
$ws.Range("J126").Value = 3033.2856
$ws.Range("K126").Value = 9249.999899999999
$ws.Range("L126").Value = 9099.856800000001
$ws.Range("M126").Value = -6779.999899999999
$ws.Range("N126").Value = -14039.8568
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 8413.25
$ws.Range("I100").Value = 13851.5
$ws.Range("J100").Value = 2975
$ws.Range("K100").Value = 13851.5
$ws.Range("L100").Value = 2975
$ws.Range("M100").Value = -13310.5
$ws.Range("N100").Value = -4057
$ws.Range("H122").Value = 22505822
$ws.Range("I122").Value = 19235920
$ws.Range("J122").Value = 28578500
$ws.Range("K122").Value = 57707760
$ws.Range("L122").Value = 85735500
$ws.Range("M122").Value = -57705310
$ws.Range("N122").Value = -85740400
$ws.Range("H136").Value = 30608424
$ws.Range("I136").Value = 41668132
$ws.Range("K136").Value = 125004396
$ws.Range("M136").Value = -125001846
